# Weekly update: prepend a new week's worth of price data (rows 186-187)
# for "Vega Monumental Concepción - Perejil", pushing the existing rows
# 186-217 down to 188-219.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right above the current row 186 (first of the two
# inserts lands at 186 and pushes the old 186 down to 187; the second
# insert at 186 again pushes that down to 188, leaving two blank rows at
# 186-187 with formatting inherited from the row that used to be there).
$ws.Rows.Item(186).Insert()
$ws.Rows.Item(186).Insert()

# --- Row 186: "Primera" quality entry for the new week ---
$ws.Range("A186").Value2 = 11
$ws.Range("B186").Value = "Vega Monumental Concepción"
$ws.Range("C186").Value = "Bíobío"
$ws.Range("D186").Value2 = 45034
$ws.Range("E186").Value2 = 8
$ws.Range("F186").Value2 = 100112044
$ws.Range("G186").Value = "Perejil"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value2 = 200
$ws.Range("K186").Value2 = 700
$ws.Range("L186").Value2 = 800
$ws.Range("M186").Value2 = 750
$ws.Range("N186").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O186").Value = "Región de Ñuble"
$ws.Range("P186").Value2 = 750
$ws.Range("Q186").Value2 = 1
$ws.Range("R186").Value = "Hortaliza"

# --- Row 187: "Segunda" quality entry for the new week ---
$ws.Range("A187").Value2 = 11
$ws.Range("B187").Value = "Vega Monumental Concepción"
$ws.Range("C187").Value = "Bíobío"
$ws.Range("D187").Value2 = 45034
$ws.Range("E187").Value2 = 8
$ws.Range("F187").Value2 = 100112044
$ws.Range("G187").Value = "Perejil"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Segunda"
$ws.Range("J187").Value2 = 100
$ws.Range("K187").Value2 = 600
$ws.Range("L187").Value2 = 600
$ws.Range("M187").Value2 = 600
$ws.Range("N187").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O187").Value = "Región de Ñuble"
$ws.Range("P187").Value2 = 600
$ws.Range("Q187").Value2 = 1
$ws.Range("R187").Value = "Hortaliza"
